# Fix the missing apostrophe in "the authors' reconstruction_minimal.py script"
# on slide 7 ("Tested the authors' 'reconstruction_minimal.py' script, but
# same result, just noise" -> "Tested the authors' ' ...").
#
# The paragraph is re-written in two steps: first to a short placeholder that
# shares no characters with the final text, then to the final text. This
# keeps the run's original formatting (<a:rPr lang="hu-HU"/>) while avoiding
# the engine splitting the paragraph into multiple runs around the single
# changed character (it diffs old vs. new text and only collapses back to a
# single run when there is no common prefix/suffix to preserve).

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(7)
$shape = $slide.Shapes.Item(3)
$textRange = $shape.TextFrame.TextRange
$paragraph = $textRange.Paragraphs(4, 1)

$paragraph.Text = "X"
$paragraph.Text = "Tested the authors’ ‘reconstruction_minimal.py’ script, but same result, just noise"
